$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.146.09"
$ws.Range("E2").Value = "  -4.39%  "

$ws.Range("D3").Value = "1.654.17"
$ws.Range("E3").Value = "  -3.25%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.37"
$ws.Range("E5").Value = "  -3.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5106"
$ws.Range("E6").Value = "  -3.50%  "

$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2588"
$ws.Range("E8").Value = "  -2.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06399"
$ws.Range("E9").Value = "  -3.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.95"
$ws.Range("E10").Value = "  -4.37%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07811"
$ws.Range("E11").Value = "  +1.67%  "

$ws.Range("D12").Value = "1.658.89"
$ws.Range("E12").Value = "  -2.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.286"
$ws.Range("E13").Value = "  -4.95%  "

$ws.Range("D14").Value = "1.884.04"
$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5528"
$ws.Range("E15").Value = "  -5.01%  "

$ws.Range("D16").Value = "0.0₅8021"
$ws.Range("E16").Value = "  -2.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.97"
$ws.Range("E17").Value = "  -5.86%  "

$ws.Range("D18").Value = "26.187.02"
$ws.Range("E18").Value = "  -4.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.68"
$ws.Range("E20").Value = "  -5.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.410"
$ws.Range("E21").Value = "  -4.84%  "

$ws.Range("E22").Value = "  -3.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.033"
$ws.Range("E23").Value = "  +0.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.26"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.726"
$ws.Range("E26").Value = "  +2.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1165"
$ws.Range("E27").Value = "  -3.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.974"
$ws.Range("E28").Value = "  -3.88%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.77"
$ws.Range("E29").Value = "  -3.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05126"
$ws.Range("E30").Value = "  -4.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.244"
$ws.Range("E31").Value = "  -3.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.352"
$ws.Range("E32").Value = "  -3.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.223"
$ws.Range("E33").Value = "  -5.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.561"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.744"
$ws.Range("E35").Value = "  -4.09%  "

$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9278"
$ws.Range("E37").Value = "  -2.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5689"
$ws.Range("E38").Value = "  -3.01%  "

$ws.Range("D39").Value = "1.153.75"
$ws.Range("E39").Value = "  +6.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01592"
$ws.Range("E40").Value = "  -2.83%  "

$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8357"
$ws.Range("E42").Value = "  -0.98%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.643"
$ws.Range("E43").Value = "  -2.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.23"
$ws.Range("E44").Value = "  -0.90%  "

$ws.Range("D45").Value = "1.793.93"
$ws.Range("E45").Value = "  -3.16%  "

$ws.Range("D46").Value = "0.0₈115"
$ws.Range("E46").Value = "  -1.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4548"
$ws.Range("E47").Value = "  +0.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.70"
$ws.Range("E48").Value = "  -3.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("E49").Value = "  -0.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.920"
$ws.Range("E50").Value = "  -2.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05042"
$ws.Range("E51").Value = "  -3.61%  "
